$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Single-cell corrections in column C (covid_deaths counts) ---
$ws.Cells.Item(762, 3).Value = 6
$ws.Cells.Item(874, 3).Value = 9
$ws.Cells.Item(941, 3).Value = 26
$ws.Cells.Item(947, 3).Value = 2
$ws.Cells.Item(956, 3).Value = 32
$ws.Cells.Item(958, 3).Value = 3
$ws.Cells.Item(960, 3).Value = 17
$ws.Cells.Item(961, 3).Value = 28
$ws.Cells.Item(962, 3).Value = 2
$ws.Cells.Item(964, 3).Value = 4
$ws.Cells.Item(967, 3).Value = 27

# --- Rows 968-983: updated/added date, age-group, covid_deaths records ---
$dateFormat = "YYYY-MM-DD HH:MM:SS"

$data = @(
    @(968, 44155, "50-59", 1),
    @(969, 44155, "60-69", 8),
    @(970, 44155, "70-79", 10),
    @(971, 44155, "80+", 23),
    @(972, 44156, "40-49", 1),
    @(973, 44156, "60-69", 9),
    @(974, 44156, "70-79", 11),
    @(975, 44156, "80+", 25),
    @(976, 44157, "30-39", 1),
    @(977, 44157, "50-59", 2),
    @(978, 44157, "60-69", 6),
    @(979, 44157, "70-79", 11),
    @(980, 44157, "80+", 27),
    @(981, 44158, "60-69", 7),
    @(982, 44158, "70-79", 11),
    @(983, 44158, "80+", 17)
)

foreach ($row in $data) {
    $r = $row[0]
    $dateVal = $row[1]
    $ageGrp = $row[2]
    $deaths = $row[3]

    $ws.Cells.Item($r, 1).Value = $dateVal
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 2).Value = $ageGrp
    $ws.Cells.Item($r, 3).Value = $deaths
}
